$d = $word.ActiveDocument

$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.End = $r1.End - 1
$r1.Text = "Uma circunferência de raio R é inscrita em um triângulo equilátero de lado L. Calcule a medida do ângulo formado pelos lados do triângulo correspondentes às cordas da circunferência."

$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.End = $r2.End - 1
$r2.Text = "Uma pirâmide é construída com base quadrada e altura de 10m. Qual é a área total da pirâmide?"
